$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the obsolete "Resolving-Mac" sending-cluster rows (10-13); new data stops at row 9
$ws.Rows("10:13").Delete()

# Refresh TPM-derived metrics for the remaining rows (2-9)
$ws.Range("I2").Value = 0.7099439172299504
$ws.Range("J2").Value = 0.7099439172299504
$ws.Range("M2").Value = 8.430598666666667
$ws.Range("N2").Value = 25.291796
$ws.Range("O2").Value = 0.173137200317126
$ws.Range("P2").Value = 0.1731372003171259
$ws.Range("Q2").Value = 3.587630021801778
$ws.Range("R2").Value = 32.288670196216
$ws.Range("S2").Value = 0.122917702211367
$ws.Range("T2").Value = 0.122917702211367
$ws.Range("I3").Value = 0.7099439172299504
$ws.Range("J3").Value = 0.7099439172299504
$ws.Range("O3").Value = 0.4685067724286191
$ws.Range("P3").Value = 0.468506772428619
$ws.Range("S3").Value = 0.3326135332667348
$ws.Range("T3").Value = 0.3326135332667348
$ws.Range("I4").Value = 0.7099439172299504
$ws.Range("J4").Value = 0.7099439172299504
$ws.Range("M4").Value = 5.125375333333333
$ws.Range("N4").Value = 15.376126
$ws.Range("O4").Value = 0.1052586145864599
$ws.Range("P4").Value = 0.1052586145864599
$ws.Range("Q4").Value = 2.181096639266222
$ws.Range("R4").Value = 19.629869753396
$ws.Range("S4").Value = 0.07472771316170895
$ws.Range("T4").Value = 0.07472771316170895
$ws.Range("I5").Value = 0.7099439172299504
$ws.Range("J5").Value = 0.7099439172299504
$ws.Range("M5").Value = 12.32411466666667
$ws.Range("N5").Value = 36.972344
$ws.Range("O5").Value = 0.2530974126677951
$ws.Range("P5").Value = 0.2530974126677951
$ws.Range("Q5").Value = 5.24451056424711
$ws.Range("R5").Value = 47.200595078224
$ws.Range("S5").Value = 0.1796849685901397
$ws.Range("T5").Value = 0.1796849685901397
$ws.Range("I6").Value = 0.2900560827700495
$ws.Range("J6").Value = 0.2900560827700495
$ws.Range("M6").Value = 8.430598666666667
$ws.Range("N6").Value = 25.291796
$ws.Range("O6").Value = 0.173137200317126
$ws.Range("P6").Value = 0.1731372003171259
$ws.Range("Q6").Value = 1.465769175982667
$ws.Range("R6").Value = 13.191922583844
$ws.Range("S6").Value = 0.05021949810575893
$ws.Range("T6").Value = 0.05021949810575892
$ws.Range("I7").Value = 0.2900560827700495
$ws.Range("J7").Value = 0.2900560827700495
$ws.Range("O7").Value = 0.4685067724286191
$ws.Range("P7").Value = 0.468506772428619
$ws.Range("S7").Value = 0.1358932391618843
$ws.Range("T7").Value = 0.1358932391618843
$ws.Range("I8").Value = 0.2900560827700495
$ws.Range("J8").Value = 0.2900560827700495
$ws.Range("M8").Value = 5.125375333333333
$ws.Range("N8").Value = 15.376126
$ws.Range("O8").Value = 0.1052586145864599
$ws.Range("P8").Value = 0.1052586145864599
$ws.Range("Q8").Value = 0.8911131315793333
$ws.Range("R8").Value = 8.020018184213999
$ws.Range("S8").Value = 0.03053090142475096
$ws.Range("T8").Value = 0.03053090142475096
$ws.Range("I9").Value = 0.2900560827700495
$ws.Range("J9").Value = 0.2900560827700495
$ws.Range("M9").Value = 12.32411466666667
$ws.Range("N9").Value = 36.972344
$ws.Range("O9").Value = 0.2530974126677951
$ws.Range("P9").Value = 0.2530974126677951
$ws.Range("Q9").Value = 2.142707548290666
$ws.Range("R9").Value = 19.284367934616
$ws.Range("S9").Value = 0.07341244407765536
$ws.Range("T9").Value = 0.07341244407765536
